# Insert two new rows at row 125 (pushing the existing rows 125-144 down
# to 127-146) and populate the two new rows with the new price records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before the current row 125.
$ws.Range("A125:A126").EntireRow.Insert()

# New row 125: Moscatel rosada
$ws.Cells.Item(125, 1).Value = 2
$ws.Cells.Item(125, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(125, 3).Value = "Coquimbo"
$ws.Cells.Item(125, 4).Value = 44644
$ws.Cells.Item(125, 5).Value = 4
$ws.Cells.Item(125, 6).Value = "Fruta"
$ws.Cells.Item(125, 7).Value = 100109
$ws.Cells.Item(125, 8).Value = "Uva"
$ws.Cells.Item(125, 9).Value = 100109001
$ws.Cells.Item(125, 10).Value = "Uva"
$ws.Cells.Item(125, 11).Value = "Moscatel rosada"
$ws.Cells.Item(125, 12).Value = "Primera"
$ws.Cells.Item(125, 13).Value = 200
$ws.Cells.Item(125, 14).Value = 13000
$ws.Cells.Item(125, 15).Value = 14000
$ws.Cells.Item(125, 16).Value = 13500
$ws.Cells.Item(125, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(125, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(125, 19).Value = 750
$ws.Cells.Item(125, 20).Value = 18

# New row 126: Red Globe
$ws.Cells.Item(126, 1).Value = 2
$ws.Cells.Item(126, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(126, 3).Value = "Coquimbo"
$ws.Cells.Item(126, 4).Value = 44644
$ws.Cells.Item(126, 5).Value = 4
$ws.Cells.Item(126, 6).Value = "Fruta"
$ws.Cells.Item(126, 7).Value = 100109
$ws.Cells.Item(126, 8).Value = "Uva"
$ws.Cells.Item(126, 9).Value = 100109001
$ws.Cells.Item(126, 10).Value = "Uva"
$ws.Cells.Item(126, 11).Value = "Red Globe"
$ws.Cells.Item(126, 12).Value = "Primera"
$ws.Cells.Item(126, 13).Value = 240
$ws.Cells.Item(126, 14).Value = 7000
$ws.Cells.Item(126, 15).Value = 8000
$ws.Cells.Item(126, 16).Value = 7500
$ws.Cells.Item(126, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(126, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(126, 19).Value = 417
$ws.Cells.Item(126, 20).Value = 18
